$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 5.910000000000001
$ws.Range("A9").Value = -21.723
$ws.Range("B12").Value = 5.927000000000001
$ws.Range("D15").Value = -8.270999999999999
$ws.Range("A18").Value = -22.156
$ws.Range("A20").Value = -20.454
$ws.Range("B26").Value = 5.761
$ws.Range("A27").Value = -21.188
$ws.Range("B27").Value = 5.695000000000001
$ws.Range("B29").Value = 5.688
$ws.Range("B37").Value = 8.376000000000001
$ws.Range("B38").Value = 5.763
$ws.Range("D38").Value = -8.669
$ws.Range("D44").Value = -7.474000000000001
$ws.Range("B51").Value = 5.14
$ws.Range("D51").Value = -8.199999999999999
$ws.Range("B55").Value = 5.705
$ws.Range("D57").Value = -8.101000000000001
$ws.Range("D63").Value = -7.336999999999999
$ws.Range("A69").Value = -21.586
$ws.Range("B69").Value = 5.688000000000001
$ws.Range("B70").Value = 5.412
$ws.Range("D70").Value = -6.797
$ws.Range("A76").Value = -20.396
$ws.Range("A82").Value = -22.23
$ws.Range("B83").Value = 5.736
$ws.Range("D99").Value = -8.138
$ws.Range("B102").Value = 7.231
